# Update the confusion-matrix derived metrics table with refreshed values
# (re-run of the underlying experiment produced new counts/metrics).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - DT
$ws.Range("B2").Value = 17045.0
$ws.Range("C2").Value = 4426.0
$ws.Range("D2").Value = 27655.0
$ws.Range("E2").Value = 5631.0
$ws.Range("F2").Value = 0.7938614875879093
$ws.Range("G2").Value = 0.7516757805609455
$ws.Range("H2").Value = 0.8620367195536299
$ws.Range("I2").Value = 0.7721929009898748
$ws.Range("M2").Value = 0.8444483996005123

# Row 3 - RF
$ws.Range("B3").Value = 14060.0
$ws.Range("C3").Value = 8517.0
$ws.Range("D3").Value = 23742.0
$ws.Range("E3").Value = 8764.0
$ws.Range("F3").Value = 0.6227576737387607
$ws.Range("G3").Value = 0.616018226428321
$ws.Range("H3").Value = 0.7359806565609597
$ws.Range("I3").Value = 0.6193696174093082
$ws.Range("M3").Value = 0.7418318726976337

# Row 4 - ZeroR
$ws.Range("B4").Value = 0.0
$ws.Range("C4").Value = 0.0
$ws.Range("D4").Value = 32259.0
$ws.Range("E4").Value = 22824.0
$ws.Range("F4").Value = ""
$ws.Range("G4").Value = 0.0
$ws.Range("H4").Value = 1.0
$ws.Range("I4").Value = ""

# Row 5 - NB
$ws.Range("B5").Value = 16149.0
$ws.Range("C5").Value = 5813.0
$ws.Range("D5").Value = 26446.0
$ws.Range("E5").Value = 6675.0
$ws.Range("F5").Value = 0.7353155450323285
$ws.Range("G5").Value = 0.7075446898002103
$ws.Range("H5").Value = 0.8198022257354537
$ws.Range("I5").Value = 0.721162863394811
$ws.Range("J5").Value = 0.7732875841911298
$ws.Range("K5").Value = 0.530275637811701
$ws.Range("L5").Value = 0.7616079775747583
$ws.Range("M5").Value = 0.8488871879585596

# Row 6 - XGB
$ws.Range("B6").Value = 16536.0
$ws.Range("C6").Value = 4366.0
$ws.Range("D6").Value = 27893.0
$ws.Range("E6").Value = 6288.0
$ws.Range("F6").Value = 0.7911204669409626
$ws.Range("G6").Value = 0.7245005257623555
$ws.Range("H6").Value = 0.8646579249201773
$ws.Range("I6").Value = 0.7563463385628688
$ws.Range("J6").Value = 0.8065827932392934
$ws.Range("K6").Value = 0.5965042006071125
$ws.Range("L6").Value = 0.7914828622334509
$ws.Range("M6").Value = 0.874645775103436

# Row 7 - k-NN
$ws.Range("B7").Value = 16702.0
$ws.Range("C7").Value = 4670.0
$ws.Range("D7").Value = 27589.0
$ws.Range("E7").Value = 6122.0
$ws.Range("F7").Value = 0.7814897997379749
$ws.Range("G7").Value = 0.7317735716789344
$ws.Range("H7").Value = 0.855234198208252
$ws.Range("I7").Value = 0.7558150058828853
